# Draft code for disagg linear regression
#
# Inserts a new "Univariate models, no STL" section at the very top of the
# document: a Heading1 title followed by a console-style (HTML Preformatted,
# dark-shaded) results block for "Linear regression (all features)", mirroring
# the layout already used further down in the document for the other models.

$d = $word.ActiveDocument

$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

# One HTML-Preformatted "console line": dark (2D2D2D) shading, light-grey
# (CCCCCC) Lucida Console text. The very first line of the block also carries
# an explicit en-US language tag on the paragraph mark (matching how the
# document's other "--------" separators are authored); the rest inherit it.
function ConsoleLineXml([string]$text, [bool]$withLang) {
    $langRun = '<w:lang w:val="en-US"/>'
    $lang = ''
    if ($withLang) { $lang = $langRun }
    return '<w:p><w:pPr><w:pStyle w:val="HTMLPreformatted"/><w:shd w:val="clear" w:color="auto" w:fill="2D2D2D"/><w:rPr><w:rFonts w:ascii="Lucida Console" w:hAnsi="Lucida Console"/><w:color w:val="CCCCCC"/>' + $lang + '</w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Lucida Console" w:hAnsi="Lucida Console"/><w:color w:val="CCCCCC"/>' + $lang + '</w:rPr><w:t>' + $text + '</w:t></w:r></w:p>'
}

$heading = '<w:p xmlns:w="' + $wNs + '"><w:pPr><w:pStyle w:val="Heading1"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Univariate models, no STL</w:t></w:r></w:p>'

$lines = @(
    (ConsoleLineXml "--------" $true),
    (ConsoleLineXml "Model = Linear regression (all features)" $false),
    (ConsoleLineXml "MAE: mean = 79.5002, sd = 199.5047, min = 0.0, max = 2476.2983" $false),
    (ConsoleLineXml "MSE: mean = 67418.3049, sd = 346389.5931, min = 0.0, max = 7386625.3219" $false),
    (ConsoleLineXml "RMSE: mean = 97.7129, sd = 240.5629, min = 0.0, max = 2717.8347" $false),
    (ConsoleLineXml "RMSLE: mean = 0.4078, sd = 0.2533, min = 0.0, max = 2.0986" $false),
    (ConsoleLineXml "--------" $false)
)

$trailingBlank = '<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p>'

$blockXml = $heading + ($lines -join '') + $trailingBlank

# Insert the whole block as raw OOXML right at the start of the body, i.e.
# immediately before the existing "Global models, no STL, 1782 series"
# heading paragraph.
$d.Range(0, 0).InsertXML($blockXml)

Write-Host "Inserted univariate linear-regression block at top of document"
